$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits right after the
#    "info@fenlogic.com." sentence near the top of the document.
#    (It gets re-created further down, around the new XBMC answer.)
# ------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
    # no-op if it doesn't exist for some reason
}

# ------------------------------------------------------------------
# 2) Merge the trailing runs of the "Ultimately yes..." paragraph
#    into a single run with the combined sentence (this also removes
#    the run-splits "R" / "unning two screens" / "SDRAM bandwidth" /
#    etc., folding everything into one <w:r>).
# ------------------------------------------------------------------
$oldText = "Ultimately yes. At the moment the drivers for two independent screens are under development. (In the demo video we cheated a bit). . Running two screens (or three with the DSI screen?) will use up a lot of SDRAM bandwidth and not all resolutions on all screens will be possible. "
$newText = "Ultimately yes. At the moment the drivers for two independent screens are under development. (In the demo video we cheated a bit). . Running two screens (or three with the DSI screen?) will use up a lot of SDRAM bandwidth and not all resolutions on all screens will be possible. "
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# Locate that paragraph (still the same one, now consolidated into one run).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Ultimately yes.*possible. `r") {
        $targetPara = $cand
        break
    }
}

# ------------------------------------------------------------------
# 3) Insert the new "Will it work for XBMC?" Q and "Yes, I have been
#    told..." A paragraphs right after it.
# ------------------------------------------------------------------
$targetPara.Range.InsertParagraphAfter()
$qPara = $targetPara.Next()

$qXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:keepNext/><w:spacing w:before="120"/><w:contextualSpacing w:val="0"/><w:rPr><w:b/><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Will it work for XBMC</w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$qPara.Range.InsertXML($qXml)

$qPara2 = $targetPara.Next()
$qPara2.Range.InsertParagraphAfter()
$aPara = $qPara2.Next()

$aXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Yes, I have been told the change are already </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>in the latest release</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$aPara.Range.InsertXML($aXml)

# ------------------------------------------------------------------
# 4) Refresh the cached PAGE field in the footer so it reflects the
#    document's real page count (the extra content above pushes the
#    footer's last-rendered page number from 1 to 6).
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$frange = $ftr.Range
$frange.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "6", 2) | Out-Null

Write-Output "done"
